$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 35 already contains the "We will update our paper very soon." annotation,
# but its politeness_score (column B) was stored as text. Convert it to a real number.
$ws.Range("B35").Value = 3

# Insert a new row of annotation data as row 36.
$ws.Range("A36").Value = "Ying Tang"

# Column B for this new row keeps its politeness_score as text "3" (not a number).
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "3"
$ws.Range("B36").Style = "Normal"

$ws.Range("C36").Value = "无"
$ws.Range("D36").Value = "DFT"
$ws.Range("E36").Value = "RES"
$ws.Range("F36").Value = "2e6daeb6-f5b1-42e4-9927-e16202e5fb2e"
$ws.Range("G36").Value = "H1cWzoxA-_annotated.xlsx"
$ws.Range("H36").Value = 'For example ,when I use the cr dataset, "python sc_main.py --network_type exp_context_fusion --context_fusion_method wblock --model_dir_suffix training --dataset_type cr --gpu 0 " the result is not the 84.48 as the paper,I could only get 84.30 after several times.'
